# -----------------------------------------------------------------------
# Refined metadata to be additional tab
#
# 1. Update the "time_taken" (column F) timestamps on the existing
#    "data" sheet (rows 2-24) to reflect the later re-run of the
#    PanelApp query.
# 2. Add a new "metadata" worksheet (placed after "data") that records
#    the panel-level metadata (name / id / version / query time / etc.)
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$dataws = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1. Refresh the per-row query timestamps on the "data" sheet (col F)
# ---------------------------------------------------------------------
$dataws.Range("F2").Value2  = "2021-10-05 14:33:17.131209"
$dataws.Range("F3").Value2  = "2021-10-05 14:33:17.131217"
$dataws.Range("F4").Value2  = "2021-10-05 14:33:17.131220"
$dataws.Range("F5").Value2  = "2021-10-05 14:33:17.131223"
$dataws.Range("F6").Value2  = "2021-10-05 14:33:17.131226"
$dataws.Range("F7").Value2  = "2021-10-05 14:33:17.131229"
$dataws.Range("F8").Value2  = "2021-10-05 14:33:17.131232"
$dataws.Range("F9").Value2  = "2021-10-05 14:33:17.131235"
$dataws.Range("F10").Value2 = "2021-10-05 14:33:17.131238"
$dataws.Range("F11").Value2 = "2021-10-05 14:33:17.131241"
$dataws.Range("F12").Value2 = "2021-10-05 14:33:17.131243"
$dataws.Range("F13").Value2 = "2021-10-05 14:33:17.131246"
$dataws.Range("F14").Value2 = "2021-10-05 14:33:17.131249"
$dataws.Range("F15").Value2 = "2021-10-05 14:33:17.131251"
$dataws.Range("F16").Value2 = "2021-10-05 14:33:17.131254"
$dataws.Range("F17").Value2 = "2021-10-05 14:33:17.131257"
$dataws.Range("F18").Value2 = "2021-10-05 14:33:17.131260"
$dataws.Range("F19").Value2 = "2021-10-05 14:33:17.131263"
$dataws.Range("F20").Value2 = "2021-10-05 14:33:17.131266"
$dataws.Range("F21").Value2 = "2021-10-05 14:33:17.131269"
$dataws.Range("F22").Value2 = "2021-10-05 14:33:17.131271"
$dataws.Range("F23").Value2 = "2021-10-05 14:33:17.131274"
$dataws.Range("F24").Value2 = "2021-10-05 14:33:17.131277"

# ---------------------------------------------------------------------
# 2. Add the new "metadata" worksheet right after "data"
# ---------------------------------------------------------------------
$metaws = $wb.Worksheets.Add($null, $dataws)
$metaws.Name = "metadata"

# --- header row (row 1), bold/bordered/centred like the "data" header ---
$dataws.Range("B1:F1").Copy()
$metaws.Range("B1:F1").PasteSpecial(-4122)   # xlPasteFormats

$metaws.Range("B1").Value2 = "data_name"
$metaws.Range("C1").Value2 = "data_id"
$metaws.Range("D1").Value2 = "data_version"
$metaws.Range("E1").Value2 = "data_version_created"
$metaws.Range("F1").Value2 = "panel_query_time"

# G1 needs the same header style as B1:F1 as well
$dataws.Range("F1").Copy()
$metaws.Range("G1").PasteSpecial(-4122)      # xlPasteFormats
$metaws.Range("G1").Value2 = "panel_get_request"

# --- data row (row 2) ---
# A2 shares the same style as the data sheet's "index" column (A2)
$dataws.Range("A2").Copy()
$metaws.Range("A2").PasteSpecial(-4122)      # xlPasteFormats
$metaws.Range("A2").Value2 = 0

$metaws.Range("B2").Value2 = "Blepharophimosis"
$metaws.Range("C2").Value2 = 55

# D2 must stay textual ("1.0"), not be coerced into the number 1.
# Build it as text on a scratch cell, then copy only the *value* (which
# carries the text cell-type) into D2, so D2 keeps the sheet's default
# (unstyled) cell format.
$scratch = $metaws.Range("ZZ100")
$scratch.NumberFormat = "@"
$scratch.Value2 = "1.0"
$scratch.Copy()
$metaws.Range("D2").PasteSpecial(-4163)      # xlPasteValues
$scratch.EntireRow.Delete()

$metaws.Range("E2").Value2 = "2021-06-06T00:31:12.359682Z"
$metaws.Range("F2").Value2 = "2021-10-05 14:33:17.127522"
$metaws.Range("G2").Value2 = "https://panelapp.agha.umccr.org/api/v1/panels/55/?format=json"

[void]$metaws.Range("A1").Select()
